# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Jenova_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1749.875
$ws.Range("J17").Value = 1749.875
$ws.Range("L17").Value = 5249.625
$ws.Range("N17").Value = -5585.625
$ws.Range("H40").Value = 9931.25
$ws.Range("J40").Value = 18372.625
$ws.Range("L40").Value = 18372.625
$ws.Range("N40").Value = -18722.625
$ws.Range("H53").Value = 23810886
$ws.Range("J53").Value = 1342.6
$ws.Range("L53").Value = 1342.6
$ws.Range("N53").Value = -2616.6
$ws.Range("H62").Value = 6254093
$ws.Range("I62").Value = 8931349
$ws.Range("J62").Value = 7162.3335
$ws.Range("K62").Value = 8931349
$ws.Range("L62").Value = 7162.3335
$ws.Range("M62").Value = -8930725
$ws.Range("N62").Value = -8410.333500000001
$ws.Range("H65").Value = 6254093
$ws.Range("I65").Value = 8931349
$ws.Range("J65").Value = 7162.3335
$ws.Range("K65").Value = 44656745
$ws.Range("L65").Value = 35811.6675
$ws.Range("M65").Value = -44653625
$ws.Range("N65").Value = -42051.6675
$ws.Range("H131").Value = 5580.8125
$ws.Range("I131").Value = 1503.6086
$ws.Range("J131").Value = 16000.333
$ws.Range("K131").Value = 4510.825800000001
$ws.Range("L131").Value = 48000.999
$ws.Range("M131").Value = 529.1741999999995
$ws.Range("N131").Value = -58080.999
$ws.Range("H133").Value = 49166
$ws.Range("J133").Value = 49166
$ws.Range("L133").Value = 49166
$ws.Range("N133").Value = -59286
$ws.Range("H135").Value = 4401.5947
$ws.Range("I135").Value = 3918.9092
$ws.Range("J135").Value = 8383.75
$ws.Range("K135").Value = 35270.1828
$ws.Range("L135").Value = 75453.75
$ws.Range("M135").Value = -32735.1828
$ws.Range("N135").Value = -80523.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3357.7415
$ws.Range("I32").Value = 3327.818
$ws.Range("K32").Value = 3327.818
$ws.Range("M32").Value = -3040.818
$ws.Range("H41").Value = 7790.8335
$ws.Range("I41").Value = 3349
$ws.Range("J41").Value = 30000
$ws.Range("K41").Value = 3349
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = -2935
$ws.Range("N41").Value = -30828
$ws.Range("H61").Value = 4635.1665
$ws.Range("I61").Value = 2952.75
$ws.Range("K61").Value = 2952.75
$ws.Range("M61").Value = -2740.75
$ws.Range("H74").Value = 852130.5
$ws.Range("I74").Value = 1002553.8
$ws.Range("K74").Value = 1002553.8
$ws.Range("M74").Value = -1001679.8
$ws.Range("H77").Value = 852130.5
$ws.Range("I77").Value = 1002553.8
$ws.Range("K77").Value = 5012769
$ws.Range("M77").Value = -5008401
$ws.Range("H102").Value = 2710.111
$ws.Range("I102").Value = 2423.875
$ws.Range("K102").Value = 2423.875
$ws.Range("M102").Value = -801.875
$ws.Range("H132").Value = 273007.8
$ws.Range("I132").Value = 328706.6
$ws.Range("J132").Value = 10427.714
$ws.Range("K132").Value = 986119.7999999999
$ws.Range("L132").Value = 31283.142
$ws.Range("M132").Value = -983589.7999999999
$ws.Range("N132").Value = -36343.142
$ws.Range("H136").Value = 4635.1665
$ws.Range("I136").Value = 2952.75
$ws.Range("K136").Value = 8858.25
$ws.Range("M136").Value = -6308.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 799.75
$ws.Range("I107").Value = 799.75
$ws.Range("K107").Value = 799.75
$ws.Range("M107").Value = 1120.25
$ws.Range("H132").Value = 50978.547
$ws.Range("J132").Value = 50978.547
$ws.Range("L132").Value = 50978.547
$ws.Range("N132").Value = -61098.547

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 432.19232
$ws.Range("I7").Value = 527.1
$ws.Range("J7").Value = 115.833336
$ws.Range("K7").Value = 527.1
$ws.Range("L7").Value = 115.833336
$ws.Range("M7").Value = -414.1
$ws.Range("N7").Value = -341.833336
$ws.Range("H20").Value = 81697.5
$ws.Range("J20").Value = 81697.5
$ws.Range("L20").Value = 81697.5
$ws.Range("N20").Value = -82169.5
$ws.Range("H30").Value = 81697.5
$ws.Range("J30").Value = 81697.5
$ws.Range("L30").Value = 81697.5
$ws.Range("N30").Value = -81879.5
$ws.Range("H51").Value = 20244.818
$ws.Range("I51").Value = 14999.167
$ws.Range("J51").Value = 26539.6
$ws.Range("K51").Value = 14999.167
$ws.Range("L51").Value = 26539.6
$ws.Range("M51").Value = -14263.167
$ws.Range("N51").Value = -28011.6
$ws.Range("H61").Value = 20244.818
$ws.Range("I61").Value = 14999.167
$ws.Range("J61").Value = 26539.6
$ws.Range("K61").Value = 14999.167
$ws.Range("L61").Value = 26539.6
$ws.Range("M61").Value = -14651.167
$ws.Range("N61").Value = -27235.6
$ws.Range("H62").Value = 3107.7778
$ws.Range("I62").Value = 2499.1667
$ws.Range("J62").Value = 4325
$ws.Range("K62").Value = 2499.1667
$ws.Range("L62").Value = 4325
$ws.Range("M62").Value = -1875.1667
$ws.Range("N62").Value = -5573
$ws.Range("H65").Value = 3107.7778
$ws.Range("I65").Value = 2499.1667
$ws.Range("J65").Value = 4325
$ws.Range("K65").Value = 12495.8335
$ws.Range("L65").Value = 21625
$ws.Range("M65").Value = -9375.833500000001
$ws.Range("N65").Value = -27865
$ws.Range("H128").Value = 81697.5
$ws.Range("J128").Value = 81697.5
$ws.Range("L128").Value = 81697.5
$ws.Range("N128").Value = -91657.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 449.33334
$ws.Range("I20").Value = 449.33334
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1348.00002
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1121.00002
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H87").Value = 5749.6665
$ws.Range("I87").Value = 5749.6665
$ws.Range("K87").Value = 17248.9995
$ws.Range("M87").Value = -16000.9995
$ws.Range("H90").Value = 5749.6665
$ws.Range("I90").Value = 5749.6665
$ws.Range("K90").Value = 51746.9985
$ws.Range("M90").Value = -45506.9985
$ws.Range("H107").Value = 20671.818
$ws.Range("I107").Value = 598
$ws.Range("J107").Value = 26273.814
$ws.Range("K107").Value = 1794
$ws.Range("L107").Value = 78821.442
$ws.Range("M107").Value = 126
$ws.Range("N107").Value = -82661.442

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 107.64
$ws.Range("I2").Value = 98.888885
$ws.Range("K2").Value = 98.888885
$ws.Range("M2").Value = 14.111115
$ws.Range("H55").Value = 16710.334
$ws.Range("J55").Value = 16724.889
$ws.Range("L55").Value = 16724.889
$ws.Range("N55").Value = -17378.889
$ws.Range("H132").Value = 24172.305
$ws.Range("I132").Value = 2322.675
$ws.Range("K132").Value = 6968.025000000001
$ws.Range("M132").Value = -4438.025000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2128.611
$ws.Range("I22").Value = 2105.4
$ws.Range("K22").Value = 2105.4
$ws.Range("M22").Value = -1810.4
$ws.Range("H27").Value = 2128.611
$ws.Range("I27").Value = 2105.4
$ws.Range("K27").Value = 2105.4
$ws.Range("M27").Value = -1998.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3333.1667
$ws.Range("I81").Value = 3333.1667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6666.3334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -5605.3334
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 3333.1667
$ws.Range("I84").Value = 3333.1667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 33331.667
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -28027.667
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 1393.2858
$ws.Range("I113").Value = 1538.9166
$ws.Range("J113").Value = 519.5
$ws.Range("K113").Value = 4616.7498
$ws.Range("L113").Value = 1558.5
$ws.Range("M113").Value = -2446.7498
$ws.Range("N113").Value = -5898.5
